$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 38463210
$ws.Range("I11").Value = 38463210
$ws.Range("K11").Value = 38463210
$ws.Range("M11").Value = -38463070
$ws.Range("H15").Value = 2998.8838
$ws.Range("I15").Value = 2998.8838
$ws.Range("K15").Value = 8996.651400000001
$ws.Range("M15").Value = -8827.651400000001
$ws.Range("H17").Value = 1480.1818
$ws.Range("J17").Value = 1480.1818
$ws.Range("L17").Value = 4440.5454
$ws.Range("N17").Value = -4776.5454
$ws.Range("H28").Value = 503.55554
$ws.Range("I28").Value = 350.55
$ws.Range("K28").Value = 350.55
$ws.Range("M28").Value = 134.45
$ws.Range("H41").Value = 2024.1765
$ws.Range("I41").Value = 1938.4546
$ws.Range("J41").Value = 2181.3333
$ws.Range("K41").Value = 1938.4546
$ws.Range("L41").Value = 2181.3333
$ws.Range("M41").Value = -1498.4546
$ws.Range("N41").Value = -3061.3333
$ws.Range("H70").Value = 2068.4736
$ws.Range("I70").Value = 2062
$ws.Range("J70").Value = 2073.182
$ws.Range("K70").Value = 6186
$ws.Range("L70").Value = 6219.545999999999
$ws.Range("M70").Value = -5916
$ws.Range("N70").Value = -6759.545999999999
$ws.Range("H73").Value = 2068.4736
$ws.Range("I73").Value = 2062
$ws.Range("J73").Value = 2073.182
$ws.Range("K73").Value = 6186
$ws.Range("L73").Value = 6219.545999999999
$ws.Range("M73").Value = -5250
$ws.Range("N73").Value = -8091.545999999999
$ws.Range("H107").Value = 240.09525
$ws.Range("I107").Value = 230.72223
$ws.Range("K107").Value = 230.72223
$ws.Range("M107").Value = 1689.27777

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14744184
$ws.Range("I61").Value = 25005162
$ws.Range("J61").Value = 85644.5
$ws.Range("K61").Value = 25005162
$ws.Range("L61").Value = 85644.5
$ws.Range("M61").Value = -25004950
$ws.Range("N61").Value = -86068.5
$ws.Range("H102").Value = 15577.9
$ws.Range("I102").Value = 18723.625
$ws.Range("J102").Value = 2995
$ws.Range("K102").Value = 18723.625
$ws.Range("L102").Value = 2995
$ws.Range("M102").Value = -17101.625
$ws.Range("N102").Value = -6239
$ws.Range("H110").Value = 1595.3846
$ws.Range("I110").Value = 1723.8182
$ws.Range("K110").Value = 1723.8182
$ws.Range("M110").Value = 321.1818000000001
$ws.Range("H132").Value = 9907.333000000001
$ws.Range("I132").Value = 3750
$ws.Range("J132").Value = 22222
$ws.Range("K132").Value = 11250
$ws.Range("L132").Value = 66666
$ws.Range("M132").Value = -8720
$ws.Range("N132").Value = -71726
$ws.Range("H136").Value = 14744184
$ws.Range("I136").Value = 25005162
$ws.Range("J136").Value = 85644.5
$ws.Range("K136").Value = 75015486
$ws.Range("L136").Value = 256933.5
$ws.Range("M136").Value = -75012936
$ws.Range("N136").Value = -262033.5
$ws.Range("H139").Value = 79493.75
$ws.Range("J139").Value = 79493.75
$ws.Range("L139").Value = 79493.75
$ws.Range("N139").Value = -89773.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8371.157999999999
$ws.Range("I99").Value = 15979.571
$ws.Range("K99").Value = 15979.571
$ws.Range("M99").Value = -14481.571
$ws.Range("H105").Value = 1847.5
$ws.Range("I105").Value = 1459.6666
$ws.Range("K105").Value = 1459.6666
$ws.Range("M105").Value = 287.3334
$ws.Range("H107").Value = 2911.5334
$ws.Range("I107").Value = 2050.7693
$ws.Range("K107").Value = 2050.7693
$ws.Range("M107").Value = -130.7692999999999
$ws.Range("H134").Value = 52763.35
$ws.Range("I134").Value = 1987.9375
$ws.Range("J134").Value = 255865
$ws.Range("K134").Value = 5963.8125
$ws.Range("L134").Value = 767595
$ws.Range("M134").Value = -3428.8125
$ws.Range("N134").Value = -772665

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -888
$ws.Range("H31").Value = 875387.9
$ws.Range("I31").Value = 28370.334
$ws.Range("J31").Value = 1117392.9
$ws.Range("K31").Value = 28370.334
$ws.Range("L31").Value = 1117392.9
$ws.Range("M31").Value = -28075.334
$ws.Range("N31").Value = -1117982.9
$ws.Range("H34").Value = 875387.9
$ws.Range("I34").Value = 28370.334
$ws.Range("J34").Value = 1117392.9
$ws.Range("K34").Value = 28370.334
$ws.Range("L34").Value = 1117392.9
$ws.Range("M34").Value = -28168.334
$ws.Range("N34").Value = -1117796.9
$ws.Range("H47").Value = 25000
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H94").Value = 3954.9412
$ws.Range("J94").Value = 4347.909
$ws.Range("L94").Value = 4347.909
$ws.Range("N94").Value = -5249.909
$ws.Range("H99").Value = 2915.3157
$ws.Range("I99").Value = 2873.0667
$ws.Range("K99").Value = 2873.0667
$ws.Range("M99").Value = -1375.0667
$ws.Range("H126").Value = 2915.3157
$ws.Range("I126").Value = 2873.0667
$ws.Range("K126").Value = 8619.2001
$ws.Range("M126").Value = -6149.2001
$ws.Range("H134").Value = 386317.97
$ws.Range("I134").Value = 436480.22
$ws.Range("K134").Value = 1309440.66
$ws.Range("M134").Value = -1306905.66

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 419.7143
$ws.Range("I14").Value = 419.7143
$ws.Range("K14").Value = 1259.1429
$ws.Range("M14").Value = -1086.1429
$ws.Range("H33").Value = 2835.3076
$ws.Range("I33").Value = 2737.1428
$ws.Range("K33").Value = 16422.8568
$ws.Range("M33").Value = -16139.8568

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 66743.836
$ws.Range("I5").Value = 80113.25
$ws.Range("J5").Value = 40005
$ws.Range("K5").Value = 80113.25
$ws.Range("L5").Value = 40005
$ws.Range("M5").Value = -80001.25
$ws.Range("N5").Value = -40229
$ws.Range("H97").Value = 2703.4443
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 2000
$ws.Range("N97").Value = -2992
$ws.Range("H122").Value = 1670.75
$ws.Range("I122").Value = 1226.8
$ws.Range("K122").Value = 3680.4
$ws.Range("M122").Value = -1230.4
$ws.Range("H126").Value = 4437.6
$ws.Range("I126").Value = 4498.25
$ws.Range("J126").Value = 4195
$ws.Range("K126").Value = 13494.75
$ws.Range("L126").Value = 12585
$ws.Range("M126").Value = -11024.75
$ws.Range("N126").Value = -17525
$ws.Range("H132").Value = 200052000
$ws.Range("I132").Value = 333366660
$ws.Range("K132").Value = 1000099980
$ws.Range("M132").Value = -1000097450
$ws.Range("H139").Value = 41983.332
$ws.Range("J139").Value = 41983.332
$ws.Range("L139").Value = 41983.332
$ws.Range("N139").Value = -52263.332

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 58807.832
$ws.Range("I7").Value = 1963.6666
$ws.Range("J7").Value = 115652
$ws.Range("K7").Value = 1963.6666
$ws.Range("L7").Value = 115652
$ws.Range("M7").Value = -1851.6666
$ws.Range("N7").Value = -115876
$ws.Range("H23").Value = 2500
$ws.Range("I23").Value = 1250
$ws.Range("K23").Value = 1250
$ws.Range("M23").Value = -1020
$ws.Range("H93").Value = 66675530
$ws.Range("J93").Value = 2113
$ws.Range("L93").Value = 2113
$ws.Range("N93").Value = -4609
$ws.Range("H126").Value = 58807.832
$ws.Range("I126").Value = 1963.6666
$ws.Range("J126").Value = 115652
$ws.Range("K126").Value = 5890.9998
$ws.Range("L126").Value = 346956
$ws.Range("M126").Value = -3420.9998
$ws.Range("N126").Value = -351896
$ws.Range("H132").Value = 33387.723
$ws.Range("I132").Value = 6977.9585
$ws.Range("J132").Value = 86207.25
$ws.Range("K132").Value = 20933.8755
$ws.Range("L132").Value = 258621.75
$ws.Range("M132").Value = -18403.8755
$ws.Range("N132").Value = -263681.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 97687.125
$ws.Range("J2").Value = 27166.334
$ws.Range("L2").Value = 27166.334
$ws.Range("N2").Value = -27390.334
$ws.Range("H41").Value = 23488.5
$ws.Range("J41").Value = 23488.5
$ws.Range("L41").Value = 23488.5
$ws.Range("N41").Value = -24268.5
$ws.Range("H62").Value = 8341334
$ws.Range("I62").Value = 8348.315000000001
$ws.Range("J62").Value = 40006680
$ws.Range("K62").Value = 8348.315000000001
$ws.Range("L62").Value = 40006680
$ws.Range("M62").Value = -7724.315000000001
$ws.Range("N62").Value = -40007928
$ws.Range("H65").Value = 8341334
$ws.Range("I65").Value = 8348.315000000001
$ws.Range("J65").Value = 40006680
$ws.Range("K65").Value = 41741.575
$ws.Range("L65").Value = 200033400
$ws.Range("M65").Value = -38621.575
$ws.Range("N65").Value = -200039640
$ws.Range("H76").Value = 44999
$ws.Range("J76").Value = 44999
$ws.Range("L76").Value = 44999
$ws.Range("N76").Value = -45629
$ws.Range("H79").Value = 44999
$ws.Range("J79").Value = 44999
$ws.Range("L79").Value = 44999
$ws.Range("N79").Value = -47183
$ws.Range("H100").Value = 1378.5217
$ws.Range("I100").Value = 1376.4762
$ws.Range("K100").Value = 2752.9524
$ws.Range("M100").Value = -2211.9524
$ws.Range("H107").Value = 100001710
$ws.Range("I107").Value = 100001710
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 300005130
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -300003210
$ws.Range("N107").ClearContents()
$ws.Range("H114").Value = 48750
$ws.Range("J114").Value = 48750
$ws.Range("L114").Value = 48750
$ws.Range("N114").Value = -57428
$ws.Range("H126").Value = 6019.2
$ws.Range("I126").Value = 5024
$ws.Range("K126").Value = 15072
$ws.Range("M126").Value = -12602
$ws.Range("H136").Value = 793.38464
$ws.Range("I136").Value = 827.5
$ws.Range("K136").Value = 2482.5
$ws.Range("M136").Value = 67.5
